$d = $word.ActiveDocument

# --- Step 1: merge the three runs of the "Magic the gathering..." paragraph
# into a single run (the visible text stays identical, the runs are just
# consolidated into one, as in the OOXML diff).
$mergedText = "Magic the gathering é um jogo de cartas feito por Richard Garfield em uma garagem em Seattle Washington com seus amigos de faculdade de matématica. a primeira coleção de magic chamada Alpha rápidamente tomou a costa Oeste por tempestade vendendo todas as cartas disponivéis nas lojas onde as cartas se encontravam, com esse sucesso estrondoso ele viu a oportunidade e lançou beta com três vezes mais cartas e acabou acontecendo outro sellout, com isso lançou unlimited com 17,700 cartas e ele começou a criar mais cartas e coleções para expandir o jogo."

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Magic the gathering*expandir o jogo.*") {
        $r = $p.Range
        [void]$r.MoveEnd(1, -1)   # exclude the paragraph mark
        # Force an actual text diff (so the runs really get consolidated)
        # by changing the text first and then setting it to the final value.
        $r.Text = "."
        $p2 = $p
        $r2 = $p2.Range
        [void]$r2.MoveEnd(1, -1)
        $r2.Text = $mergedText
        break
    }
}

# --- Step 2: insert two new paragraphs right after the "Hoje em dia..."
# paragraph, before the following (empty, size 32) paragraph.
$para1 = "O jogo tem vários formatos porém o que irei falar é sobre o standard que é um modo de jogo que dois jogadores vão montar um deck de 60 cartas entre elas cartas das coleções mais novas, e seu objetivo é utilizar essas cartas para diminuir a vida do seu oponente para zero"
$para2 = "Eu conheci magic quando fui para uma loja de cartas perto de casa com o meu irmão e o mesmo me ensinou a jogar; jogo esse que já adorei logo de cara vendo todas as possibilidades de cartas e jeitos diferentes de jogar, meu primeiro deck foi um preto e branco ganho de vida onde a ideia era simplesmente ter vida o suficiente para o meu oponente não conseguir me vencer"

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Hoje em dia magic*mostrados nas cartas*") {
        $anchor = $p.Range
        [void]$anchor.Collapse(0)   # wdCollapseEnd
        [void]$anchor.InsertParagraphAfter()

        $newPara1 = $p.Next()
        $newPara1.Range.Text = $para1

        $anchor2 = $newPara1.Range
        [void]$anchor2.Collapse(0)
        [void]$anchor2.InsertParagraphAfter()

        $newPara2 = $newPara1.Next()
        $newPara2.Range.Text = $para2
        break
    }
}
